$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.284.38"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.631.11"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.05"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.38"
$ws.Range("E6").Value = "  -3.01%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.571"
$ws.Range("E8").Value = "  -4.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.639.71"
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("E10").Value = "  -5.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.106"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.086.55"
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.258.26"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.27"
$ws.Range("E16").Value = "  -2.06%  "
$ws.Range("E17").Value = "  -1.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.631.65"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.64"
$ws.Range("E19").Value = "  -2.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "341.80"
$ws.Range("E20").Value = "  -2.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.46"
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.14"
$ws.Range("E22").Value = "  -1.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.996"
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.16"
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("E25").Value = "  -1.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("E27").Value = "  -2.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0811"
$ws.Range("E28").Value = "  -3.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.07"
$ws.Range("E29").Value = "  -3.60%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.99"
$ws.Range("E32").Value = "  -4.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.00"
$ws.Range("E33").Value = "  -2.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.64"
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("E35").Value = "  -5.78%  "
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("E37").Value = "  -4.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.865"
$ws.Range("E38").Value = "  +2.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.62"
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("E40").Value = "  -4.22%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.65"
$ws.Range("E41").Value = "  -3.73%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "291.84"
$ws.Range("E42").Value = "  +1.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.628"
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("E46").Value = "  -2.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.48"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("E48").Value = "  +0.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0232"
$ws.Range("E49").Value = "  -2.25%  "
$ws.Range("E50").Value = "  -5.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.964.18"
$ws.Range("E51").Value = "  +0.02%  "
